# Insert a new weekly price record as row 106 in the daily-logic subset
# sheet, pushing the existing rows 106-156 down to 107-157 (the sheet's
# dimension grows from A1:R156 to A1:R157).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 106:156 down by one row, copying formatting (date style on
# column D, etc.) from the row being pushed down - this matches Excel's
# native "Insert Row" behaviour.
$ws.Rows(106).Insert()

# Populate the newly-inserted row 106 with the new observation.
$ws.Range("A106").Value = 11
$ws.Range("B106").Value = "Vega Monumental Concepción"
$ws.Range("C106").Value = "Bíobío"
$ws.Range("D106").Value = 44846
$ws.Range("E106").Value = 8
$ws.Range("F106").Value = 100112032
$ws.Range("G106").Value = "Zapallo italiano"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 220
$ws.Range("K106").Value = 18000
$ws.Range("L106").Value = 20000
$ws.Range("M106").Value = 19091
$ws.Range("N106").Value = "$/caja 50 unidades"
$ws.Range("O106").Value = "Región de O'Higgins"
$ws.Range("P106").Value = 382
$ws.Range("Q106").Value = 50
$ws.Range("R106").Value = "Hortaliza"
